# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    3  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    4  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    5  = @(0.1554434735375247, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 2.461472421528573)
    6  = @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.705647867635037)
    7  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    8  = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    9  = @(3.182878228561681, 9.226618575922256, 0.7127328510149897, 6.48142807727062, 19.60365773276954)
    10 = @(0.3464964993005633, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 1.051601690082842)
    11 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    12 = @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 5.964442013611383)
    13 = @(0.1554434735375247, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 2.461472421528573)
    14 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 3.594575437922795)
    15 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    16 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    17 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    18 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    19 = @(0.7287194209349384, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 4.64879039072104)
    20 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    21 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    22 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 22.31973251085698)
    23 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    24 = @(1.505614041169197, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 20.64246832346449)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
